$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the X_UTM / Y_UTM header labels in B1 / C1
$ws.Range("B1").Value = "Y_UTM"
$ws.Range("C1").Value = "X_UTM"

# Add new "area" column header in AF1, matching the header style used by
# the other header cells (bold, bordered, centered) by copying AE1's format
$ws.Range("AF1").Value = "area"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)

# Fill AF2:AF477 with the value 4 for every data row
for ($r = 2; $r -le 477; $r++) {
    $ws.Cells.Item($r, 32).Value = 4
}
